$d = $word.ActiveDocument

# 1. Insert two new paragraphs at the very start of the document body:
#    - a Heading2 styled paragraph with the "Adding Text..." sentence
#    - a plain paragraph with "Jonah Ripley is responsible..." (with the
#      grammar-check proofErr markers around "changes" that Word's own
#      editor would have produced)
$introXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:pPr><w:pStyle w:val="Heading2"/><w:pBdr>' +
            '<w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/>' +
            '<w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr></w:pPr>' +
            '<w:r><w:t>Adding Text to test GIT commits and pushes.</w:t></w:r>' +
            '</w:p>' +
            '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
            '<w:r><w:t xml:space="preserve">Jonah Ripley is responsible for these document </w:t></w:r>' +
            '<w:proofErr w:type="gramStart"/>' +
            '<w:r><w:t>changes</w:t></w:r>' +
            '<w:proofErr w:type="gramEnd"/>' +
            '</w:p>'

$startRange = $d.Range(0, 0)
$startRange.InsertXML($introXml)

# 2. Fill in the two previously-empty cells in the "Team Member" /
#    "Roles and Responsibilities" table row.
$table = $d.Tables(1)
$table.Cell(2, 1).Range.Text = "Jonah Ripley"
$table.Cell(2, 2).Range.Text = "Testing GIT functions for in-class activities"
